# Apply updated odds values for Jogos_da_Semana_FlashScore_2025-03-13.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.32
$ws.Range("H2").Value = 2.75
$ws.Range("J2").Value = 3
$ws.Range("L2").Value = 4.2
$ws.Range("P2").Value = 2.3
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.44
$ws.Range("S2").Value = 4.7
$ws.Range("T2").Value = 1.15
$ws.Range("Y2").Value = 5.8
$ws.Range("AC2").Value = 23
$ws.Range("AD2").Value = 45
$ws.Range("AF2").Value = 5.5
$ws.Range("AH2").Value = 120
$ws.Range("AN2").Value = 40
$ws.Range("G3").Value = 1.9
$ws.Range("I3").Value = 3.35
$ws.Range("J3").Value = 2.42
$ws.Range("K3").Value = 2.27
$ws.Range("L3").Value = 3.7
$ws.Range("O3").Value = 1.2
$ws.Range("P3").Value = 3.65
$ws.Range("R3").Value = 2.05
$ws.Range("S3").Value = 2.42
$ws.Range("T3").Value = 1.44
$ws.Range("X3").Value = 2.1
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 10.25
$ws.Range("AB3").Value = 16.5
$ws.Range("AC3").Value = 14
$ws.Range("AD3").Value = 22
$ws.Range("AF3").Value = 7.5
$ws.Range("AJ3").Value = 12.5
$ws.Range("AK3").Value = 20
$ws.Range("AL3").Value = 11.75
$ws.Range("AM3").Value = 45
$ws.Range("AN3").Value = 27
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 3.8
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.2
$ws.Range("K4").Value = 2.25
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 8.5
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.98
$ws.Range("S4").Value = 3
$ws.Range("T4").Value = 1.36
$ws.Range("U4").Value = 1.36
$ws.Range("V4").Value = 3
$ws.Range("W4").Value = 1.8
$ws.Range("X4").Value = 1.91
$ws.Range("Y4").Value = 7.5
$ws.Range("Z4").Value = 8
$ws.Range("AB4").Value = 12
$ws.Range("AD4").Value = 26
$ws.Range("AE4").Value = 11
$ws.Range("AF4").Value = 7.5
$ws.Range("AH4").Value = 51
$ws.Range("AJ4").Value = 15
$ws.Range("AL4").Value = 17
$ws.Range("AN4").Value = 41
$ws.Range("AO4").Value = 41
$ws.Range("AP4").Value = 1.98
$ws.Range("AQ4").Value = 1.83
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 12
$ws.Range("O5").Value = 1.13
$ws.Range("P5").Value = 5.5
$ws.Range("G6").Value = 1.22
$ws.Range("H6").Value = 6.5
$ws.Range("I6").Value = 11
$ws.Range("J6").Value = 1.57
$ws.Range("L6").Value = 8.5
$ws.Range("M6").Value = 19
$ws.Range("N6").Value = 1.03
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 1.73
$ws.Range("U6").Value = 1.22
$ws.Range("V6").Value = 4
$ws.Range("Y6").Value = 10
$ws.Range("AF6").Value = 13
$ws.Range("AI6").Value = 700
$ws.Range("AK6").Value = 51
$ws.Range("AL6").Value = 29
$ws.Range("AM6").Value = 126
$ws.Range("AN6").Value = 67
$ws.Range("G7").Value = 2.32
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 3.2
$ws.Range("J7").Value = 2.87
$ws.Range("K7").Value = 2.02
$ws.Range("L7").Value = 3.7
$ws.Range("M7").Value = 1.09
$ws.Range("N7").Value = 6.2
$ws.Range("O7").Value = 1.38
$ws.Range("P7").Value = 2.8
$ws.Range("Q7").Value = 2.12
$ws.Range("R7").Value = 1.65
$ws.Range("S7").Value = 3.6
$ws.Range("T7").Value = 1.24
$ws.Range("U7").Value = 1.42
$ws.Range("V7").Value = 2.65
$ws.Range("W7").Value = 1.82
$ws.Range("X7").Value = 1.88
$ws.Range("Y7").Value = 7
$ws.Range("Z7").Value = 10.75
$ws.Range("AA7").Value = 9
$ws.Range("AB7").Value = 24
$ws.Range("AC7").Value = 20
$ws.Range("AD7").Value = 32
$ws.Range("AE7").Value = 6.2
$ws.Range("AF7").Value = 5.8
$ws.Range("AG7").Value = 14
$ws.Range("AH7").Value = 70
$ws.Range("AJ7").Value = 8.75
$ws.Range("AK7").Value = 16.5
$ws.Range("AL7").Value = 11
$ws.Range("AM7").Value = 45
$ws.Range("AN7").Value = 30
$ws.Range("AO7").Value = 40

Write-Output "Applied changes"
